$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("eye")

# Disable combinations that are difficult to see:
# - Column P (Circuit sclera_00014.jpg) is disabled for every Metal Iris row (2-10)
# - Row 7 (Metal Iris_00151.jpg) is also disabled for column Q (Circuit sclera_00015.jpg)
$ws.Range("P2:P10").ClearContents()
$ws.Range("Q7").ClearContents()

# Update the view state to match: scroll right and select Q8
$ws.Range("Q8").Select()
$excel.ActiveWindow.ScrollColumn = 7
